# Weekly driver report update for 2025-04-20
# Updates the "Bad Drivers" and "Good Drivers" tables on the
# "Driver Summary" sheet with the latest weekly roaming data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Bad Drivers (rows 3-7) + Totals (row 8)
# ---------------------------------------------------------------

# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 22.110.1.1 (unchanged name/count)
$ws.Range("C3").Value = 414
$ws.Range("D3").Value = 98.09999999999999

# Row 4 now reports the 22.240.0.6 driver
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.240.0.6"
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 1472
$ws.Range("D4").Value = 98.5

# Row 5: Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.3.2 (unchanged name/count)
$ws.Range("C5").Value = 635

# Row 6: Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2 (unchanged name/count)
$ws.Range("C6").Value = 122
$ws.Range("D6").Value = 98.8

# Row 7 now reports the 23.40.0.4 driver
$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4"
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 614
$ws.Range("D7").Value = 98.90000000000001

# Totals row
$ws.Range("C8").Value = 3257

# ---------------------------------------------------------------
# Good Drivers (Roaming > 99.8%) - rows 16-21
# ---------------------------------------------------------------

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B16").Value = 445055
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").Value = "'2024-11-10"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B17").Value = 77849
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").Value = "'2021-08-18"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B18").Value = 34244
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "'2021-04-27"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B19").Value = 59673
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "'2020-08-05"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B20").Value = 113652
# D20 (100) is unchanged
$ws.Range("E20").Value = "'2020-01-06"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B21").Value = 56018
# D21 (100) and E21 (2019-12-14) are unchanged, left as-is
